$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.592.83"
$ws.Range("E2").Value = "  -0.18%  "
$ws.Range("D3").Value = "2.272.07"
$ws.Range("E3").Value = "  -0.31%  "
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("D5").Value = "'119.59"
$ws.Range("E5").Value = "  +5.12%  "
$ws.Range("D6").Value = "'266.03"
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("D7").Value = "'0.646"
$ws.Range("E7").Value = "  +3.61%  "
$ws.Range("E8").Value = "  +0.21%  "
$ws.Range("D9").Value = "'0.623"
$ws.Range("E9").Value = "  +2.13%  "
$ws.Range("D10").Value = "'47.71"
$ws.Range("E10").Value = "  -1.10%  "
$ws.Range("E11").Value = "  +1.05%  "
$ws.Range("D12").Value = "'9.23"
$ws.Range("E12").Value = "  +4.75%  "
$ws.Range("E13").Value = "  -1.19%  "
$ws.Range("D14").Value = "'15.47"
$ws.Range("E14").Value = "  -1.09%  "
$ws.Range("E15").Value = "  +4.45%  "
$ws.Range("D16").Value = "2.613.16"
$ws.Range("E16").Value = "  -0.26%  "
$ws.Range("D17").Value = "2.261.43"
$ws.Range("E17").Value = "  -0.82%  "
$ws.Range("D18").Value = "43.538.18"
$ws.Range("E18").Value = "  +0.21%  "
$ws.Range("E19").Value = "  +1.56%  "
$ws.Range("D20").Value = "'6.90"
$ws.Range("E20").Value = "  -1.38%  "
$ws.Range("D21").Value = "'72.14"
$ws.Range("E21").Value = "  +0.42%  "
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("D23").Value = "'235.80"
$ws.Range("E23").Value = "  +1.59%  "
$ws.Range("E24").Value = "  -3.78%  "
$ws.Range("E25").Value = "  +0.60%  "
$ws.Range("D26").Value = "'12.01"
$ws.Range("E26").Value = "  +4.70%  "
$ws.Range("D28").Value = "'41.79"
$ws.Range("E28").Value = "  +2.00%  "
$ws.Range("D29").Value = "'3.39"
$ws.Range("E29").Value = "  -0.29%  "
$ws.Range("E30").Value = "  -0.06%  "
$ws.Range("D31").Value = "'172.01"
$ws.Range("E31").Value = "  -0.65%  "
$ws.Range("D32").Value = "'21.62"
$ws.Range("E32").Value = "  +0.94%  "
$ws.Range("D33").Value = "'0.0918"
$ws.Range("E33").Value = "  +0.66%  "
$ws.Range("D34").Value = "'5.76"
$ws.Range("E34").Value = "  +2.35%  "
$ws.Range("E35").Value = "  +3.07%  "
$ws.Range("D36").Value = "'0.0386"
$ws.Range("E36").Value = "  +10.60%  "
$ws.Range("D37").Value = "'4.21"
$ws.Range("E37").Value = "  +12.68%  "
$ws.Range("E38").Value = "  -0.61%  "
$ws.Range("E39").Value = "  +1.53%  "
$ws.Range("E40").Value = "  +6.04%  "
$ws.Range("D41").Value = "'73.80"
$ws.Range("E41").Value = "  -0.66%  "
$ws.Range("D42").Value = "'13.76"
$ws.Range("E42").Value = "  -3.86%  "
$ws.Range("D43").Value = "'0.238"
$ws.Range("E43").Value = "  +0.85%  "
$ws.Range("D44").Value = "'0.999"
$ws.Range("E44").Value = "  -0.36%  "
$ws.Range("D45").Value = "'1.38"
$ws.Range("E45").Value = "  -0.14%  "
$ws.Range("E46").Value = "  -6.45%  "
$ws.Range("D47").Value = "'73.77"
$ws.Range("E47").Value = "  +41.21%  "
$ws.Range("E48").Value = "  +1.64%  "
$ws.Range("D49").Value = "'8.54"
$ws.Range("E49").Value = "  -1.30%  "
$ws.Range("E50").Value = "  +0.56%  "
$ws.Range("D51").Value = "'101.89"
$ws.Range("E51").Value = "  +0.39%  "
